# "Add files via upload" — adds a new bug-tracker row (A9) with a Jira link,
# gives it a medium box border + wrap text, and updates the current
# selection/scroll position to reflect where the user ended up after typing
# it in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: A9 gets the Jira/Atlassian issue link (this also appends
# the string to the shared-string table, bumping sharedStrings count/uniqueCount).
$ws.Range("A9").Value = "https://vladm2329.atlassian.net/browse/RRRR-1"

# Wrap the long URL text, matching the other description cells in the sheet.
$ws.Range("A9").WrapText = $true

# Box the new cell with a medium border on all four sides (new border +
# cellXfs entry get created automatically).
$ws.Range("A9").Borders.Weight = -4138  # xlMedium

# Leave the selection where it ended up after entering the new row
# (one row below the last filled cell), matching the saved view state.
[void]$ws.Range("B10").Select()
